$d = $word.ActiveDocument

# Locate the "GWP definition " paragraph and the empty ListParagraph
# paragraph that immediately follows it (the last two paragraphs of the
# document body, just before the sectPr).
$count = $d.Paragraphs.Count
$pGwp = $d.Paragraphs($count - 1)
$pEmpty = $d.Paragraphs($count)

# Range spanning from the start of the "GWP definition " paragraph
# through the end of the trailing empty paragraph (i.e. up to, but not
# including, the sectPr).
$full = $d.Range($pGwp.Range.Start, $pEmpty.Range.End)

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
       '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
       '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
       '<pkg:xmlData>' +
       '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
       '<w:body>' +
       '<w:p>' +
       '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr>' +
       '<w:r><w:t xml:space="preserve">GWP </w:t></w:r>' +
       '<w:r><w:t>definition</w:t></w:r>' +
       '</w:p>' +
       '<w:p/><w:p/><w:p/><w:p/>' +
       '</w:body>' +
       '</w:document>' +
       '</pkg:xmlData></pkg:part></pkg:package>'

$full.InsertXML($xml)
